$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.375.36"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3
$ws.Range("D3").Value = "2.604.71"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.94"
$ws.Range("E5").Value = "  -1.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.40"
$ws.Range("E6").Value = "  -2.70%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.92%  "

# Row 9
$ws.Range("D9").Value = "2.603.68"
$ws.Range("E9").Value = "  -0.36%  "

# Row 10
$ws.Range("E10").Value = "  +1.81%  "

# Row 11
$ws.Range("E11").Value = "  +0.02%  "

# Row 12
$ws.Range("E12").Value = "  -1.48%  "

# Row 13
$ws.Range("E13").Value = "  -2.99%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.25"
$ws.Range("E14").Value = "  -2.81%  "

# Row 15
$ws.Range("D15").Value = "3.079.33"
$ws.Range("E15").Value = "  -0.28%  "

# Row 16
$ws.Range("E16").Value = "  -2.87%  "

# Row 17
$ws.Range("D17").Value = "67.189.30"
$ws.Range("E17").Value = "  -0.57%  "

# Row 18
$ws.Range("D18").Value = "2.607.59"
$ws.Range("E18").Value = "  -0.33%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "370.20"
$ws.Range("E19").Value = "  +1.16%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.02"
$ws.Range("E20").Value = "  -2.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.35"
$ws.Range("E21").Value = "  -3.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  -2.80%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.75"
$ws.Range("E23").Value = "  -4.31%  "

# Row 24
$ws.Range("E24").Value = "  -3.41%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.21"
$ws.Range("E25").Value = "  +4.63%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.91"
$ws.Range("E27").Value = "  -2.18%  "

# Row 28
$ws.Range("D28").Value = "2.732.13"
$ws.Range("E28").Value = "  -0.51%  "

# Row 29
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "578.01"
$ws.Range("E30").Value = "  -0.99%  "

# Row 31
$ws.Range("E31").Value = "  -6.01%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.36"
$ws.Range("E32").Value = "  -5.35%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.66"
$ws.Range("E33").Value = "  -3.53%  "

# Row 34
$ws.Range("E34").Value = "  -3.16%  "

# Row 35
$ws.Range("E35").Value = "  +0.05%  "

# Row 36
$ws.Range("E36").Value = "  -4.10%  "

# Row 37
$ws.Range("E37").Value = "  -2.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.46"
$ws.Range("E38").Value = "  +2.14%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.04"
$ws.Range("E39").Value = "  -1.86%  "

# Row 40
$ws.Range("E40").Value = "  +0.67%  "

# Row 41
$ws.Range("E41").Value = "  -1.90%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.21"
$ws.Range("E42").Value = "  -3.47%  "

# Row 43
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.12"
$ws.Range("E43").Value = "  +4.17%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  -4.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.95"
$ws.Range("E46").Value = "  -2.30%  "

# Row 47
$ws.Range("D47").Value = "0.0₆0281"
$ws.Range("E47").Value = "  -1.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.65"
$ws.Range("E48").Value = "  -2.98%  "

# Row 49
$ws.Range("E49").Value = "  -1.52%  "

# Row 50
$ws.Range("E50").Value = "  -4.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.30"
$ws.Range("E51").Value = "  +1.33%  "
